$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 64
$ws1.Range("F4").Value = 0
$ws1.Range("F6").Value = 143
$ws1.Range("F7").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("F10").Value = 5359
$ws1.Range("F11").Value = 0
$ws1.Range("F12").Value = 6038
$ws1.Range("F14").Value = 384
$ws1.Range("F15").Value = 394
$ws1.Range("F16").Value = 30
$ws1.Range("F17").Value = 551
$ws1.Range("F18").Value = 0
$ws1.Range("F19").Value = 268
$ws1.Range("F20").Value = 142
$ws1.Range("F21").Value = 201
$ws1.Range("F23").Value = 102
$ws1.Range("F24").Value = 0
$ws1.Range("F25").Value = 1870
$ws1.Range("F26").Value = 1705
$ws1.Range("F27").Value = 47
$ws1.Range("F29").Value = 0
$ws1.Range("F31").Value = 85
$ws1.Range("F32").Value = 154
$ws1.Range("F34").Value = 2018
$ws1.Range("F35").Value = 302
$ws1.Range("F36").Value = 0
$ws1.Range("F38").Value = 5135
$ws1.Range("F41").Value = 633
$ws1.Range("F42").Value = 100
$ws1.Range("F43").Value = 166
$ws1.Range("F45").Value = 1069
$ws1.Range("F47").Value = 0
$ws1.Range("F48").Value = 59

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 3
$ws2.Range("F8").Value = 0
$ws2.Range("F15").Value = 0
$ws2.Range("F16").Value = 0
$ws2.Range("F17").Value = 7
$ws2.Range("F18").Value = 0
$ws2.Range("F19").Value = 0

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F3").Value = 64
$ws4.Range("F7").Value = 143
$ws4.Range("F8").Value = 234
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 5
$ws4.Range("F13").Value = 0
$ws4.Range("F15").Value = 6038
$ws4.Range("F16").Value = 1094
$ws4.Range("F18").Value = 394
$ws4.Range("F20").Value = 551
$ws4.Range("F22").Value = 268
$ws4.Range("F23").Value = 0
$ws4.Range("F24").Value = 201
$ws4.Range("F25").Value = 0
$ws4.Range("F27").Value = 190
$ws4.Range("F28").Value = 9899
$ws4.Range("F30").Value = 1705
$ws4.Range("F31").Value = 0
$ws4.Range("F32").Value = 2009
$ws4.Range("F33").Value = 77
$ws4.Range("F34").Value = 85
$ws4.Range("F36").Value = 7
$ws4.Range("F37").Value = 2018
$ws4.Range("F38").Value = 302
$ws4.Range("F40").Value = 5135
$ws4.Range("F42").Value = 0
$ws4.Range("F44").Value = 166
$ws4.Range("F45").Value = 1098
$ws4.Range("F48").Value = 0
$ws4.Range("F49").Value = 59
